$d = $word.ActiveDocument

# Update the cover-page team member names (Thai + English) to the new
# student's name, per the commit's content change.

$d.Content.Find.Execute("ชลกันต์", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "ชินพัฒน์", 2)

$d.Content.Find.Execute("บังเกิด", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "ลิ้มประธาน", 2)

$d.Content.Find.Execute("CHONLAGUN", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "CHINNAPAT", 2)

$d.Content.Find.Execute("BANGKERT", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "LIMPRATHAN", 2)
